$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-13
# from 2023-10-25 (serial 45224) to 2023-11-03 (serial 45233)
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = (Get-Date -Year 2023 -Month 11 -Day 3 -Hour 0 -Minute 0 -Second 0).Date
}
